$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column M ("MAR" bound column); this shifts
#    the existing M/N data to N/O and everything right of M (T:AE) to (U:AF).
$ws.Range("M1").EntireColumn.Insert()

# 2. New column header in M10 -> shared string "MAR"
$ws.Range("M10").Value = "MAR"

# 3. Fill the new M column for the existing data rows (15-25) with 0
$ws.Range("M15:M25").Value = 0

# 4. Add six new rows (26-31) below the existing table, inheriting the
#    formatting of the row above (same as the existing rows 19-25 pattern).
$ws.Range("A26:A31").EntireRow.Insert()

# Row 26: ERSOLPV5N - 2015
$ws.Range("D26").Value = 2015
$ws.Range("E26").Value = "UP"
$ws.Range("F26").Value = "CAP_BND"
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Style = "Normal 42"
$ws.Range("N26").Value = "ERSOLPV5N"

# Row 27: ERSOLPV5N - 2020
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = "UP"
$ws.Range("F27").Value = "CAP_BND"
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Style = "Normal"
$ws.Range("N27").Formula = "=N26"

# Row 28: ERSOLPV5N - 2025
$ws.Range("D28").Value = 2025
$ws.Range("E28").Value = "UP"
$ws.Range("F28").Value = "CAP_BND"
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Style = "Normal"
$ws.Range("N28").Formula = "=N27"

# Row 29: TB_H2_MAR_DKW_01 - 2015
$ws.Range("D29").Value = 2015
$ws.Range("E29").Value = "UP"
$ws.Range("F29").Value = "CAP_BND"
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Style = "Normal"
$ws.Range("N29").Value = "TB_H2_MAR_DKW_01"

# Row 30: TB_H2_MAR_DKW_01 - 2020
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = "UP"
$ws.Range("F30").Value = "CAP_BND"
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Style = "Normal"
$ws.Range("N30").Formula = "=N29"

# Row 31: TB_H2_MAR_DKW_01 - 2025
$ws.Range("D31").Value = 2025
$ws.Range("E31").Value = "UP"
$ws.Range("F31").Value = "CAP_BND"
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Style = "Normal"
$ws.Range("N31").Formula = "=N30"

$ws.Range("M10").Select()
